# Auto-generated Excel COM-interop edit script
# Applies numeric corrections to several rows across the ALC, ARM, CRP, CUL, GSM, LTW, WVR sheets
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 6469.971
$ws.Range("I137").Value = 3905.3
$ws.Range("J137").Value = 9889.532999999999
$ws.Range("K137").Value = 11715.9
$ws.Range("L137").Value = 29668.599
$ws.Range("M137").Value = -9165.900000000001
$ws.Range("N137").Value = -34768.599

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H6").Value = 18467.666
$ws.Range("I6").Value = 17700
$ws.Range("K6").Value = 17700
$ws.Range("M6").Value = -17527

$ws.Range("H45").Value = 1717.0667
$ws.Range("I45").Value = 1362.3636
$ws.Range("K45").Value = 1362.3636
$ws.Range("M45").Value = -985.3635999999999

$ws.Range("H61").Value = 5037.6597
$ws.Range("I61").Value = 5341.276
$ws.Range("J61").Value = 4548.5
$ws.Range("K61").Value = 5341.276
$ws.Range("L61").Value = 4548.5
$ws.Range("M61").Value = -5129.276
$ws.Range("N61").Value = -4972.5

$ws.Range("H74").Value = 8035.875
$ws.Range("I74").Value = 10759.4
$ws.Range("K74").Value = 10759.4
$ws.Range("M74").Value = -9885.4

$ws.Range("H77").Value = 8035.875
$ws.Range("I77").Value = 10759.4
$ws.Range("K77").Value = 53797
$ws.Range("M77").Value = -49429

$ws.Range("H109").Value = 51666.668
$ws.Range("J109").Value = 51666.668
$ws.Range("L109").Value = 51666.668
$ws.Range("N109").Value = -54440.668

$ws.Range("H136").Value = 5037.6597
$ws.Range("I136").Value = 5341.276
$ws.Range("J136").Value = 4548.5
$ws.Range("K136").Value = 16023.828
$ws.Range("L136").Value = 13645.5
$ws.Range("M136").Value = -13473.828
$ws.Range("N136").Value = -18745.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 364037.06
$ws.Range("I22").Value = 725735.75
$ws.Range("K22").Value = 725735.75
$ws.Range("M22").Value = -725385.75

$ws.Range("H31").Value = 4592.0415
$ws.Range("I31").Value = 1772
$ws.Range("K31").Value = 1772
$ws.Range("M31").Value = -1477

$ws.Range("H34").Value = 4592.0415
$ws.Range("I34").Value = 1772
$ws.Range("K34").Value = 1772
$ws.Range("M34").Value = -1570

$ws.Range("H97").Value = 0
$ws.Range("J97").Value = 0
$ws.Range("L97").Value = 0
$ws.Range("N97").ClearContents()

$ws.Range("H107").Value = 800
$ws.Range("I107").Value = 800
$ws.Range("K107").Value = 800
$ws.Range("M107").Value = 1120

$ws.Range("H132").Value = 4819.5386
$ws.Range("I132").Value = 3701.8572
$ws.Range("J132").Value = 14599.25
$ws.Range("K132").Value = 11105.5716
$ws.Range("L132").Value = 43797.75
$ws.Range("M132").Value = -8575.571599999999
$ws.Range("N132").Value = -48857.75

$ws.Range("H134").Value = 62517970
$ws.Range("I134").Value = 90922776
$ws.Range("J134").Value = 27391.8
$ws.Range("K134").Value = 272768328
$ws.Range("L134").Value = 82175.39999999999
$ws.Range("M134").Value = -272765793
$ws.Range("N134").Value = -87245.39999999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 6404639
$ws.Range("I4").Value = 7143221
$ws.Range("K4").Value = 21429663
$ws.Range("M4").Value = -21429551

$ws.Range("H39").Value = 4760.2
$ws.Range("J39").Value = 4760.2
$ws.Range("L39").Value = 14280.6
$ws.Range("N39").Value = -14868.6

$ws.Range("H50").Value = 3084
$ws.Range("I50").Value = 2700.8
$ws.Range("K50").Value = 8102.400000000001
$ws.Range("M50").Value = -7621.400000000001

$ws.Range("H53").Value = 3084
$ws.Range("I53").Value = 2700.8
$ws.Range("K53").Value = 8102.400000000001
$ws.Range("M53").Value = -7621.400000000001

$ws.Range("H88").Value = 24166
$ws.Range("J88").Value = 24166
$ws.Range("L88").Value = 72498
$ws.Range("N88").Value = -73354

$ws.Range("H91").Value = 24166
$ws.Range("J91").Value = 24166
$ws.Range("L91").Value = 72498
$ws.Range("N91").Value = -75462

$ws.Range("H107").Value = 3487.9756
$ws.Range("I107").Value = 412.85715
$ws.Range("J107").Value = 4121.0884
$ws.Range("K107").Value = 1238.57145
$ws.Range("L107").Value = 12363.2652
$ws.Range("M107").Value = 681.4285500000001
$ws.Range("N107").Value = -16203.2652

$ws.Range("H122").Value = 122409.36
$ws.Range("I122").Value = 271.2
$ws.Range("K122").Value = 2440.8
$ws.Range("M122").Value = 9.200000000000273

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H58").Value = 8767.75
$ws.Range("I58").Value = 8767.75
$ws.Range("K58").Value = 8767.75
$ws.Range("M58").Value = -8490.75

$ws.Range("H70").Value = 13356.936
$ws.Range("I70").Value = 9386.526
$ws.Range("J70").Value = 19643.416
$ws.Range("K70").Value = 9386.526
$ws.Range("L70").Value = 19643.416
$ws.Range("M70").Value = -9116.526
$ws.Range("N70").Value = -20183.416

$ws.Range("H73").Value = 13356.936
$ws.Range("I73").Value = 9386.526
$ws.Range("J73").Value = 19643.416
$ws.Range("K73").Value = 9386.526
$ws.Range("L73").Value = 19643.416
$ws.Range("M73").Value = -8450.526
$ws.Range("N73").Value = -21515.416

$ws.Range("H107").Value = 1332.7
$ws.Range("I107").Value = 453.5
$ws.Range("K107").Value = 453.5
$ws.Range("M107").Value = 1466.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 3837.186
$ws.Range("I61").Value = 2547.3242
$ws.Range("K61").Value = 2547.3242
$ws.Range("M61").Value = -2345.3242

$ws.Range("H113").Value = 3837.186
$ws.Range("I113").Value = 2547.3242
$ws.Range("K113").Value = 2547.3242
$ws.Range("M113").Value = -377.3242

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H75").Value = 44949
$ws.Range("J75").Value = 44949
$ws.Range("L75").Value = 44949
$ws.Range("N75").Value = -46821

$ws.Range("H78").Value = 44949
$ws.Range("J78").Value = 44949
$ws.Range("L78").Value = 134847
$ws.Range("N78").Value = -144207

$ws.Range("H126").Value = 7285.1333
$ws.Range("I126").Value = 5523.3335
$ws.Range("J126").Value = 14332.333
$ws.Range("K126").Value = 16570.0005
$ws.Range("L126").Value = 42996.999
$ws.Range("M126").Value = -14100.0005
$ws.Range("N126").Value = -47936.999

$ws.Range("H132").Value = 5366.479
$ws.Range("I132").Value = 3833.1865
$ws.Range("K132").Value = 11499.5595
$ws.Range("M132").Value = -8969.559499999999
